$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.425.85'
$ws.Range("E2").Value = '  +5.72%  '

$ws.Range("D3").Value = '2.601.01'
$ws.Range("E3").Value = '  +7.63%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = "'508.19"
$ws.Range("E5").Value = '  +4.05%  '

$ws.Range("D6").Value = "'156.27"
$ws.Range("E6").Value = '  +0.83%  '

$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = '  -4.21%  '

$ws.Range("D9").Value = '2.626.94'
$ws.Range("E9").Value = '  +7.89%  '

$ws.Range("E10").Value = '  +4.38%  '

$ws.Range("E11").Value = '  +3.89%  '

$ws.Range("D12").Value = "'0.344"
$ws.Range("E12").Value = '  +3.85%  '

$ws.Range("E13").Value = '  +0.98%  '

$ws.Range("D14").Value = '3.020.88'
$ws.Range("E14").Value = '  +6.83%  '

$ws.Range("D15").Value = '60.415.00'
$ws.Range("E15").Value = '  +5.64%  '

$ws.Range("D16").Value = "'21.81"
$ws.Range("E16").Value = '  +5.32%  '

$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = '  +5.18%  '

$ws.Range("D18").Value = '2.608.17'
$ws.Range("E18").Value = '  +7.31%  '

$ws.Range("D19").Value = "'4.81"
$ws.Range("E19").Value = '  +2.87%  '

$ws.Range("D20").Value = "'349.32"
$ws.Range("E20").Value = '  +8.88%  '

$ws.Range("D21").Value = "'10.43"
$ws.Range("E21").Value = '  +4.20%  '

$ws.Range("D22").Value = "'6.12"
$ws.Range("E22").Value = '  +3.52%  '

$ws.Range("E23").Value = '  +0.31%  '

$ws.Range("D24").Value = "'60.32"
$ws.Range("E24").Value = '  +3.90%  '

$ws.Range("D25").Value = "'0.423"
$ws.Range("E25").Value = '  +5.21%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = '  +3.64%  '

$ws.Range("B27").Value = 'WrappedeETH'
$ws.Range("C27").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D27").Value = '2.686.84'
$ws.Range("E27").Value = '  +6.56%  '

$ws.Range("D28").Value = "'0.982"
$ws.Range("E28").Value = '  -0.81%  '

$ws.Range("D29").Value = '0.0₃0863'
$ws.Range("E29").Value = '  +9.29%  '

$ws.Range("D30").Value = "'7.49"
$ws.Range("E30").Value = '  +2.64%  '

$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").Value = "'156.30"
$ws.Range("E32").Value = '  +3.71%  '

$ws.Range("D33").Value = "'19.33"
$ws.Range("E33").Value = '  +3.32%  '

$ws.Range("D34").Value = "'1.57"
$ws.Range("E34").Value = '  +2.86%  '

$ws.Range("D35").Value = "'5.77"
$ws.Range("E35").Value = '  +8.40%  '

$ws.Range("D36").Value = "'4.02"
$ws.Range("E36").Value = '  +6.66%  '

$ws.Range("D37").Value = "'1.21"
$ws.Range("E37").Value = '  +6.05%  '

$ws.Range("D38").Value = "'0.864"
$ws.Range("E38").Value = '  +27.35%  '

$ws.Range("D39").Value = "'0.855"
$ws.Range("E39").Value = '  +4.62%  '

$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = '  +7.16%  '

$ws.Range("D41").Value = "'3.78"
$ws.Range("E41").Value = '  +6.85%  '

$ws.Range("D42").Value = "'302.43"
$ws.Range("E42").Value = '  +8.73%  '

$ws.Range("D43").Value = "'35.70"
$ws.Range("E43").Value = '  +3.97%  '

$ws.Range("D44").Value = "'0.0571"
$ws.Range("E44").Value = '  +6.67%  '

$ws.Range("D45").Value = "'0.621"
$ws.Range("E45").Value = '  +4.27%  '

$ws.Range("D46").Value = "'0.100"
$ws.Range("E46").Value = '  -0.15%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'19.98"
$ws.Range("E47").Value = '  +11.75%  '

$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").Value = "'0.993"
$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").Value = "'4.99"
$ws.Range("E49").Value = '  +8.05%  '

$ws.Range("D50").Value = "'0.0236"
$ws.Range("E50").Value = '  +3.08%  '

$ws.Range("D51").Value = '2.051.35'
$ws.Range("E51").Value = '  +8.63%  '

# Restore default (General, no quote-prefix) formatting on cells forced to text
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
